$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '69.727.56'
$ws.Range("E2").Value2 = '  +0.67%  '
Set-TextValue $ws.Range("D3") '3.428.99'
$ws.Range("E3").Value2 = '  +1.08%  '
$ws.Range("E4").Value2 = '  +0.12%  '
Set-TextValue $ws.Range("D5") '584.02'
$ws.Range("E5").Value2 = '  -0.67%  '
$ws.Range("E6").Value2 = '  -2.09%  '
Set-TextValue $ws.Range("D7") '3.422.86'
$ws.Range("E7").Value2 = '  +1.07%  '
$ws.Range("E8").Value2 = '  +0.04%  '
Set-TextValue $ws.Range("D9") '0.591'
$ws.Range("E9").Value2 = '  -0.79%  '
$ws.Range("E10").Value2 = '  +2.46%  '
Set-TextValue $ws.Range("D11") '0.584'
$ws.Range("E11").Value2 = '  -1.03%  '
Set-TextValue $ws.Range("D12") '49.15'
$ws.Range("E12").Value2 = '  +0.84%  '
Set-TextValue $ws.Range("D13") '0.0000283'
$ws.Range("E13").Value2 = '  +0.21%  '
Set-TextValue $ws.Range("D14") '691.04'
$ws.Range("E14").Value2 = '  +1.69%  '
Set-TextValue $ws.Range("D15") '3.980.03'
$ws.Range("E15").Value2 = '  +0.96%  '
$ws.Range("E16").Value2 = '  +0.04%  '
Set-TextValue $ws.Range("D17") '69.825.37'
$ws.Range("E17").Value2 = '  +0.67%  '
Set-TextValue $ws.Range("D20") '17.69'
$ws.Range("E20").Value2 = '  -0.12%  '
Set-TextValue $ws.Range("D21") '11.40'
$ws.Range("E21").Value2 = '  +0.29%  '
$ws.Range("E22").Value2 = '  -0.53%  '
Set-TextValue $ws.Range("D23") '5.45'
$ws.Range("E23").Value2 = '  +0.71%  '
Set-TextValue $ws.Range("D24") '16.93'
$ws.Range("E24").Value2 = '  -1.16%  '
Set-TextValue $ws.Range("D25") '100.84'
$ws.Range("E25").Value2 = '  -2.68%  '
$ws.Range("E26").Value2 = '  -0.15%  '
$ws.Range("E27").Value2 = '  -2.63%  '
Set-TextValue $ws.Range("D28") '9.64'
$ws.Range("E28").Value2 = '  +0.21%  '
Set-TextValue $ws.Range("D29") '33.53'
$ws.Range("E29").Value2 = '  -1.89%  '
Set-TextValue $ws.Range("D30") '8.77'
$ws.Range("E30").Value2 = '  +0.65%  '
Set-TextValue $ws.Range("D31") '7.15'
$ws.Range("E31").Value2 = '  +2.15%  '
Set-TextValue $ws.Range("D32") '574.90'
$ws.Range("E32").Value2 = '  +3.67%  '
$ws.Range("E33").Value2 = '  -1.73%  '
Set-TextValue $ws.Range("D34") '11.02'
$ws.Range("E34").Value2 = '  -1.61%  '
Set-TextValue $ws.Range("D35") '58.21'
$ws.Range("E35").Value2 = '  +0.51%  '
$ws.Range("E36").Value2 = '  -2.57%  '
Set-TextValue $ws.Range("D38") '3.573.82'
$ws.Range("E38").Value2 = '  -3.52%  '
$ws.Range("E39").Value2 = '  -0.41%  '
Set-TextValue $ws.Range("D40") '35.31'
$ws.Range("E40").Value2 = '  +0.16%  '
Set-TextValue $ws.Range("D41") '0.0₃0738'
$ws.Range("E41").Value2 = '  +4.58%  '
$ws.Range("E42").Value2 = '  +0.65%  '
Set-TextValue $ws.Range("D43") '2.68'
$ws.Range("E43").Value2 = '  +0.24%  '
$ws.Range("E44").Value2 = '  +3.54%  '
Set-TextValue $ws.Range("D48") '2.66'
$ws.Range("E48").Value2 = '  -0.16%  '
$ws.Range("E49").Value2 = '  -0.95%  '
$ws.Range("E50").Value2 = '  -0.31%  '
Set-TextValue $ws.Range("D51") '132.39'
$ws.Range("E51").Value2 = '  +0.67%  '

$ws.Range("B18").Value2 = 'WrappedEther'
$ws.Range("C18").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D18") '3.427.57'
$ws.Range("E18").Value2 = '  +1.00%  '
$ws.Range("B19").Value2 = 'TRON'
$ws.Range("C19").Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D19") '0.122'
$ws.Range("E19").Value2 = '  +1.13%  '
$ws.Range("B45").Value2 = 'TheGraph'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D45") '0.334'
$ws.Range("E45").Value2 = '  -1.59%  '
$ws.Range("B46").Value2 = 'VeChain'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D46") '0.0420'
$ws.Range("E46").Value2 = '  -0.82%  '
